# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-22) is re-sorted from
# descending period order (1903..1809) to ascending period order
# (1809..1903). Only the values travel to their new rows - the existing
# per-row formatting (the bottom row keeps its distinct border style) is
# left exactly where it is, so we overwrite the cell values directly
# instead of using Range.Sort (which would also relocate the formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1809", "1810", "1811", "1812", "1901", "1902", "1903")
$valores  = @(31249, 31249, 31249, 31249, 31249, 31249, 26041)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value2 = $periodos[$i]
    $ws.Range("F$row").Value2 = $valores[$i]
}
